# Add three new character styles (GaNStyle, GaNParagraph, GaNLinks) and
# apply them to the runs that were newly styled per the commit diff.

$d = $word.ActiveDocument

# --- Define the new character styles -------------------------------------

$GaNStyle = $d.Styles.Add("GaNStyle", 2)
$GaNStyle.Font.Name = "Calibri"
$GaNStyle.Font.Size = 14

$GaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$GaNParagraph.Font.Name = "Calibri"
$GaNParagraph.Font.Size = 10

$GaNLinks = $d.Styles.Add("GaNLinks", 2)
$GaNLinks.Font.Name = "Calibri"
$GaNLinks.Font.Bold = $true
$GaNLinks.Font.Color = 8388608    # BGR 0x800000 -> renders as RRGGBB 000080 (navy)
$GaNLinks.Font.Size = 9.5
$GaNLinks.Font.Underline = 1      # wdUnderlineSingle

# --- Helper: apply a character style to every exact-text match -----------

function Apply-GaNStyle($searchText, $styleName) {
    $range = $d.Content
    while ($range.Find.Execute($searchText, $true, $false, $false, $false,
                                $false, $true, 1, $false, "", 0)) {
        $range.Style = $styleName
        $range.Collapse(0)
        $range.End = $d.Content.End
    }
}

# "Dates de la campanya 2022 ..." heading, repeated 4 times in the doc
Apply-GaNStyle "Dates de la campanya 2022 en què usem la  Constel·lació de Bessons 14-23 de febrer, 14-24 de març" "GaNStyle"

# "Esteu participant en una campanya mundial ..." intro paragraph
Apply-GaNStyle "Esteu participant en una campanya mundial per observar i anotar la brillantor de les estrelles més febles que es poden veure, com a mitjà per mesurar la contaminació lumínica en un lloc determinat. Localitzant i observant la  Constel·lació de Bessons a la nit i comparant la brillantor de les estrelles del cel amb la brillantor que indiquen els mapes, gent de tot el món aprendran com els llums de la seva zona contribueixen a augmentar la contaminació lumínica. Les vostres aportacions a la base de dades activa faran palesa la visibilitat del cel nocturn." "GaNParagraph"

# "Jenik Hollan, CzechGlobe (...)" credit line
Apply-GaNStyle "Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)." "GaNLinks"
